# Auto-generated edit script: update market-price-derived profit columns (H-N)
# across ALC/ARM/BSM/CRP/GSM sheets per scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 964.8333
$ws.Range("I18").Value = 867.8
$ws.Range("K18").Value = 867.8
$ws.Range("M18").Value = -583.8

$ws.Range("H132").Value = 2049.451
$ws.Range("I132").Value = 1623.1945
$ws.Range("J132").Value = 3072.4666
$ws.Range("K132").Value = 4869.583500000001
$ws.Range("L132").Value = 9217.399800000001
$ws.Range("M132").Value = -2339.583500000001
$ws.Range("N132").Value = -14277.3998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 458665.3
$ws.Range("I61").Value = 315850.78
$ws.Range("J61").Value = 913075.25
$ws.Range("K61").Value = 315850.78
$ws.Range("L61").Value = 913075.25
$ws.Range("M61").Value = -315638.78
$ws.Range("N61").Value = -913499.25

$ws.Range("H63").Value = 2658.2104
$ws.Range("I63").Value = 1642.8572
$ws.Range("J63").Value = 3250.5
$ws.Range("K63").Value = 1642.8572
$ws.Range("L63").Value = 3250.5
$ws.Range("M63").Value = -956.8571999999999
$ws.Range("N63").Value = -4622.5

$ws.Range("H66").Value = 2658.2104
$ws.Range("I66").Value = 1642.8572
$ws.Range("J66").Value = 3250.5
$ws.Range("K66").Value = 8214.286
$ws.Range("L66").Value = 16252.5
$ws.Range("M66").Value = -4782.286
$ws.Range("N66").Value = -23116.5

$ws.Range("H136").Value = 458665.3
$ws.Range("I136").Value = 315850.78
$ws.Range("J136").Value = 913075.25
$ws.Range("K136").Value = 947552.3400000001
$ws.Range("L136").Value = 2739225.75
$ws.Range("M136").Value = -945002.3400000001
$ws.Range("N136").Value = -2744325.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 9275.75
$ws.Range("I75").Value = 4413.625
$ws.Range("J75").Value = 19000
$ws.Range("K75").Value = 4413.625
$ws.Range("L75").Value = 19000
$ws.Range("M75").Value = -3477.625
$ws.Range("N75").Value = -20872

$ws.Range("H78").Value = 9275.75
$ws.Range("I78").Value = 4413.625
$ws.Range("J78").Value = 19000
$ws.Range("K78").Value = 13240.875
$ws.Range("L78").Value = 57000
$ws.Range("M78").Value = -8560.875
$ws.Range("N78").Value = -66360

$ws.Range("H99").Value = 3361.1667
$ws.Range("I99").Value = 4284.533
$ws.Range("K99").Value = 4284.533
$ws.Range("M99").Value = -2786.533

$ws.Range("H105").Value = 1686.7727
$ws.Range("I105").Value = 1553.9333
$ws.Range("K105").Value = 1553.9333
$ws.Range("M105").Value = 193.0667000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 463.6
$ws.Range("I8").Value = 463.6
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 463.6
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -323.6
$ws.Range("N8").ClearContents()

$ws.Range("H31").Value = 2400.946
$ws.Range("I31").Value = 1632.4182
$ws.Range("K31").Value = 1632.4182
$ws.Range("M31").Value = -1337.4182

$ws.Range("H34").Value = 2400.946
$ws.Range("I34").Value = 1632.4182
$ws.Range("K34").Value = 1632.4182
$ws.Range("M34").Value = -1430.4182

$ws.Range("H69").Value = 7333.3335
$ws.Range("I69").Value = 7333.3335
$ws.Range("K69").Value = 7333.3335
$ws.Range("M69").Value = -6584.3335

$ws.Range("H72").Value = 7333.3335
$ws.Range("I72").Value = 7333.3335
$ws.Range("K72").Value = 22000.0005
$ws.Range("M72").Value = -18256.0005

$ws.Range("H87").Value = 184665
$ws.Range("J87").Value = 184665
$ws.Range("L87").Value = 184665
$ws.Range("N87").Value = -187037

$ws.Range("H90").Value = 184665
$ws.Range("J90").Value = 184665
$ws.Range("L90").Value = 553995
$ws.Range("N90").Value = -565851

$ws.Range("H99").Value = 55167
$ws.Range("I99").Value = 112827.11
$ws.Range("J99").Value = 3272.9
$ws.Range("K99").Value = 112827.11
$ws.Range("L99").Value = 3272.9
$ws.Range("M99").Value = -111329.11
$ws.Range("N99").Value = -6268.9

$ws.Range("H122").Value = 1339.4286
$ws.Range("I122").Value = 858.7778
$ws.Range("J122").Value = 2204.6
$ws.Range("K122").Value = 2576.3334
$ws.Range("L122").Value = 6613.799999999999
$ws.Range("M122").Value = -126.3334
$ws.Range("N122").Value = -11513.8

$ws.Range("H126").Value = 55167
$ws.Range("I126").Value = 112827.11
$ws.Range("J126").Value = 3272.9
$ws.Range("K126").Value = 338481.33
$ws.Range("L126").Value = 9818.700000000001
$ws.Range("M126").Value = -336011.33
$ws.Range("N126").Value = -14758.7

$ws.Range("H132").Value = 1978.766
$ws.Range("I132").Value = 1088.3438
$ws.Range("J132").Value = 3878.3333
$ws.Range("K132").Value = 3265.0314
$ws.Range("L132").Value = 11634.9999
$ws.Range("M132").Value = -735.0314000000003
$ws.Range("N132").Value = -16694.9999

$ws.Range("H134").Value = 1432.4054
$ws.Range("I134").Value = 971.75
$ws.Range("K134").Value = 2915.25
$ws.Range("M134").Value = -380.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4342.5166
$ws.Range("I70").Value = 4225.6313
$ws.Range("J70").Value = 4544.409
$ws.Range("K70").Value = 4225.6313
$ws.Range("L70").Value = 4544.409
$ws.Range("M70").Value = -3955.6313
$ws.Range("N70").Value = -5084.409

$ws.Range("H73").Value = 4342.5166
$ws.Range("I73").Value = 4225.6313
$ws.Range("J73").Value = 4544.409
$ws.Range("K73").Value = 4225.6313
$ws.Range("L73").Value = 4544.409
$ws.Range("M73").Value = -3289.6313
$ws.Range("N73").Value = -6416.409

$ws.Range("H80").Value = 10290.5
$ws.Range("I80").Value = 13817.5
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 13817.5
$ws.Range("L80").Value = 5000
$ws.Range("M80").Value = -12819.5
$ws.Range("N80").Value = -6996

$ws.Range("H83").Value = 10290.5
$ws.Range("I83").Value = 13817.5
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 69087.5
$ws.Range("L83").Value = 25000
$ws.Range("M83").Value = -64095.5
$ws.Range("N83").Value = -34984

$ws.Range("H97").Value = 1520.9
$ws.Range("I97").Value = 1652.375
$ws.Range("J97").Value = 995
$ws.Range("K97").Value = 1652.375
$ws.Range("L97").Value = 995
$ws.Range("M97").Value = -1156.375
$ws.Range("N97").Value = -1987

$ws.Range("H126").Value = 2987
$ws.Range("I126").Value = 2693.36
$ws.Range("J126").Value = 3373.3684
$ws.Range("K126").Value = 8080.08
$ws.Range("L126").Value = 10120.1052
$ws.Range("M126").Value = -5610.08
$ws.Range("N126").Value = -15060.1052
